# Scheduled data-refresh: update cached Universalis market-price / profit
# columns (H:N = currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ,
# LeveProfitHQ) for the affected leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1057.6383
$ws.Cells.Item(132, 9).Value = 730
$ws.Cells.Item(132, 11).Value = 2190
$ws.Cells.Item(132, 13).Value = 340

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2132628.2
$ws.Cells.Item(138, 9).Value = 3225.5715
$ws.Cells.Item(138, 10).Value = 3036011
$ws.Cells.Item(138, 11).Value = 9676.7145
$ws.Cells.Item(138, 12).Value = 9108033
$ws.Cells.Item(138, 13).Value = -4536.7145
$ws.Cells.Item(138, 14).Value = -9118313

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 48902.09
$ws.Cells.Item(74, 9).Value = 73596.5
$ws.Cells.Item(74, 10).Value = 5686.875
$ws.Cells.Item(74, 11).Value = 73596.5
$ws.Cells.Item(74, 12).Value = 5686.875
$ws.Cells.Item(74, 13).Value = -72722.5
$ws.Cells.Item(74, 14).Value = -7434.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 48902.09
$ws.Cells.Item(77, 9).Value = 73596.5
$ws.Cells.Item(77, 10).Value = 5686.875
$ws.Cells.Item(77, 11).Value = 367982.5
$ws.Cells.Item(77, 12).Value = 28434.375
$ws.Cells.Item(77, 13).Value = -363614.5
$ws.Cells.Item(77, 14).Value = -37170.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 8336965
$ws.Cells.Item(97, 9).Value = 4124.75
$ws.Cells.Item(97, 10).Value = 13892192
$ws.Cells.Item(97, 11).Value = 4124.75
$ws.Cells.Item(97, 12).Value = 13892192
$ws.Cells.Item(97, 13).Value = -3628.75
$ws.Cells.Item(97, 14).Value = -13893184

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 13894195
$ws.Cells.Item(110, 9).Value = 6548.263
$ws.Cells.Item(110, 11).Value = 6548.263
$ws.Cells.Item(110, 13).Value = -4503.263

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 8949.875
$ws.Cells.Item(132, 9).Value = 7335.3
$ws.Cells.Item(132, 10).Value = 11640.833
$ws.Cells.Item(132, 11).Value = 22005.9
$ws.Cells.Item(132, 12).Value = 34922.499
$ws.Cells.Item(132, 13).Value = -19475.9
$ws.Cells.Item(132, 14).Value = -39982.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 100233.11
$ws.Cells.Item(133, 10).Value = 100233.11
$ws.Cells.Item(133, 12).Value = 100233.11
$ws.Cells.Item(133, 14).Value = -105293.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2524.375
$ws.Cells.Item(105, 9).Value = 1604.6857
$ws.Cells.Item(105, 11).Value = 1604.6857
$ws.Cells.Item(105, 13).Value = 142.3143

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4468881.5
$ws.Cells.Item(134, 9).Value = 6581016
$ws.Cells.Item(134, 10).Value = 9931.723
$ws.Cells.Item(134, 11).Value = 19743048
$ws.Cells.Item(134, 12).Value = 29795.169
$ws.Cells.Item(134, 13).Value = -19740513
$ws.Cells.Item(134, 14).Value = -34865.169

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6593.213
$ws.Cells.Item(31, 10).Value = 11228
$ws.Cells.Item(31, 12).Value = 11228
$ws.Cells.Item(31, 14).Value = -11818

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 6593.213
$ws.Cells.Item(34, 10).Value = 11228
$ws.Cells.Item(34, 12).Value = 11228
$ws.Cells.Item(34, 14).Value = -11632

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4633.9414
$ws.Cells.Item(58, 9).Value = 2150
$ws.Cells.Item(58, 11).Value = 2150
$ws.Cells.Item(58, 13).Value = -1947

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 6491.636
$ws.Cells.Item(99, 9).Value = 3998
$ws.Cells.Item(99, 10).Value = 6741
$ws.Cells.Item(99, 11).Value = 3998
$ws.Cells.Item(99, 12).Value = 6741
$ws.Cells.Item(99, 13).Value = -2500
$ws.Cells.Item(99, 14).Value = -9737

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4029.4
$ws.Cells.Item(122, 9).Value = 3162.5
$ws.Cells.Item(122, 10).Value = 4607.3335
$ws.Cells.Item(122, 11).Value = 9487.5
$ws.Cells.Item(122, 12).Value = 13822.0005
$ws.Cells.Item(122, 13).Value = -7037.5
$ws.Cells.Item(122, 14).Value = -18722.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 6491.636
$ws.Cells.Item(126, 9).Value = 3998
$ws.Cells.Item(126, 10).Value = 6741
$ws.Cells.Item(126, 11).Value = 11994
$ws.Cells.Item(126, 12).Value = 20223
$ws.Cells.Item(126, 13).Value = -9524
$ws.Cells.Item(126, 14).Value = -25163

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4701.636
$ws.Cells.Item(132, 9).Value = 1977.7727
$ws.Cells.Item(132, 10).Value = 7425.5
$ws.Cells.Item(132, 11).Value = 5933.3181
$ws.Cells.Item(132, 12).Value = 22276.5
$ws.Cells.Item(132, 13).Value = -3403.3181
$ws.Cells.Item(132, 14).Value = -27336.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 5042.909
$ws.Cells.Item(134, 9).Value = 1865.3334
$ws.Cells.Item(134, 10).Value = 7944.174
$ws.Cells.Item(134, 11).Value = 5596.0002
$ws.Cells.Item(134, 12).Value = 23832.522
$ws.Cells.Item(134, 13).Value = -3061.0002
$ws.Cells.Item(134, 14).Value = -28902.522

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 4633.9414
$ws.Cells.Item(136, 9).Value = 2150
$ws.Cells.Item(136, 11).Value = 6450
$ws.Cells.Item(136, 13).Value = -3900

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 857
$ws.Cells.Item(25, 9).Value = 966.5
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 2899.5
$ws.Cells.Item(25, 12).Value = 600
$ws.Cells.Item(25, 13).Value = -2730.5
$ws.Cells.Item(25, 14).Value = -938

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(30, 8).Value = 857
$ws.Cells.Item(30, 9).Value = 966.5
$ws.Cells.Item(30, 10).Value = 200
$ws.Cells.Item(30, 11).Value = 2899.5
$ws.Cells.Item(30, 12).Value = 600
$ws.Cells.Item(30, 13).Value = -2797.5
$ws.Cells.Item(30, 14).Value = -804

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 3536667.2
$ws.Cells.Item(122, 10).Value = 1001.6667
$ws.Cells.Item(122, 12).Value = 9015.0003
$ws.Cells.Item(122, 14).Value = -13915.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 6000
$ws.Cells.Item(19, 9).Value = 6000
$ws.Cells.Item(19, 11).Value = 6000
$ws.Cells.Item(19, 13).Value = -5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2515.4285
$ws.Cells.Item(80, 9).Value = 2551.6667
$ws.Cells.Item(80, 11).Value = 2551.6667
$ws.Cells.Item(80, 13).Value = -1553.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2515.4285
$ws.Cells.Item(83, 9).Value = 2551.6667
$ws.Cells.Item(83, 11).Value = 12758.3335
$ws.Cells.Item(83, 13).Value = -7766.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1026.4736
$ws.Cells.Item(97, 9).Value = 934.2
$ws.Cells.Item(97, 10).Value = 1372.5
$ws.Cells.Item(97, 11).Value = 934.2
$ws.Cells.Item(97, 12).Value = 1372.5
$ws.Cells.Item(97, 13).Value = -438.2
$ws.Cells.Item(97, 14).Value = -2364.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(118, 8).Value = 31750
$ws.Cells.Item(118, 10).Value = 31750
$ws.Cells.Item(118, 12).Value = 31750
$ws.Cells.Item(118, 14).Value = -35064

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4037881
$ws.Cells.Item(122, 9).Value = 5589989
$ws.Cells.Item(122, 10).Value = 2400
$ws.Cells.Item(122, 11).Value = 16769967
$ws.Cells.Item(122, 12).Value = 7200
$ws.Cells.Item(122, 13).Value = -16767517
$ws.Cells.Item(122, 14).Value = -12100

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3938.7942
$ws.Cells.Item(132, 9).Value = 1541.3914
$ws.Cells.Item(132, 10).Value = 8951.546
$ws.Cells.Item(132, 11).Value = 4624.174199999999
$ws.Cells.Item(132, 12).Value = 26854.638
$ws.Cells.Item(132, 13).Value = -2094.174199999999
$ws.Cells.Item(132, 14).Value = -31914.638

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4715.3887
$ws.Cells.Item(122, 9).Value = 2989.0908
$ws.Cells.Item(122, 10).Value = 7428.143
$ws.Cells.Item(122, 11).Value = 8967.2724
$ws.Cells.Item(122, 12).Value = 22284.429
$ws.Cells.Item(122, 13).Value = -6517.2724
$ws.Cells.Item(122, 14).Value = -27184.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(128, 8).Value = 49649.332
$ws.Cells.Item(128, 10).Value = 49649.332
$ws.Cells.Item(128, 12).Value = 49649.332
$ws.Cells.Item(128, 14).Value = -59609.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6261.553
$ws.Cells.Item(132, 9).Value = 2899.7827
$ws.Cells.Item(132, 10).Value = 9483.25
$ws.Cells.Item(132, 11).Value = 8699.348100000001
$ws.Cells.Item(132, 12).Value = 28449.75
$ws.Cells.Item(132, 13).Value = -6169.348100000001
$ws.Cells.Item(132, 14).Value = -33509.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 13380.897
$ws.Cells.Item(136, 10).Value = 21361.227
$ws.Cells.Item(136, 12).Value = 64083.681
$ws.Cells.Item(136, 14).Value = -69183.681

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 836.3137
$ws.Cells.Item(113, 9).Value = 748.1667
$ws.Cells.Item(113, 11).Value = 2244.5001
$ws.Cells.Item(113, 13).Value = -74.5001000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 339223.5
$ws.Cells.Item(122, 9).Value = 670747.5
$ws.Cells.Item(122, 11).Value = 2012242.5
$ws.Cells.Item(122, 13).Value = -2009792.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5280.125
$ws.Cells.Item(132, 9).Value = 6104.737
$ws.Cells.Item(132, 10).Value = 4074.923
$ws.Cells.Item(132, 11).Value = 18314.211
$ws.Cells.Item(132, 12).Value = 12224.769
$ws.Cells.Item(132, 13).Value = -15784.211
$ws.Cells.Item(132, 14).Value = -17284.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 440798.3
$ws.Cells.Item(136, 9).Value = 1833.8334
$ws.Cells.Item(136, 11).Value = 5501.5002
$ws.Cells.Item(136, 13).Value = -2951.5002
